$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Script text in A6: the trailing "{0}" placeholder line now gets a "p: " speaker prefix
$ws.Range("A6").Value = "`n@set ""game_chapter = 1""`n@showUI MapButtonUI`n@showUI QuestLogUI`n`np: {0}`n"
$ws.Rows.Item(6).AutoFit()

# 2. Column B no longer needs its own wider width - bring it back down to the sheet's
#    regular column width (8.67) so it reads the same as every other column.
$ws.Range("B:B").ColumnWidth = 7.78

# 3. The active selection moves to C6
$ws.Range("C6").Select()
